# PacBio Proposal.docx edit
#   - Title paragraph replaced with "Genomic Approach to Conservation of
#     Crayfish Biodiversity" (plain, no bold/indent formatting)
#   - Three new paragraphs inserted (Biodiversity / blank / PI: Carla Hurt)
#   - The scratch/notes block near the end of the document is removed

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Replace the bold "Using HiFi sequencing..." paragraph (4th paragraph)
#    with a single plain-text title paragraph.
# ---------------------------------------------------------------------
$titlePara = $d.Paragraphs.Item(4)
$titleStart = $titlePara.Range.Start
$titleEnd = $titlePara.Range.End
$d.Range($titleStart, $titleEnd).Delete()

$newTitleAnchor = $d.Paragraphs.Item(4)
$insRange = $d.Range($newTitleAnchor.Range.Start, $newTitleAnchor.Range.Start)
$insRange.InsertParagraphBefore()

$titleParaNow = $d.Paragraphs.Item(4)
$titleParaNow.Range.Text = "Genomic Approach to Conservation of Crayfish Biodiversity"

# ---------------------------------------------------------------------
# 2) Insert three new paragraphs ("Biodiversity", blank, "PI: Carla Hurt ")
#    between the two blank paragraphs that used to follow the title
#    paragraph directly.
# ---------------------------------------------------------------------
$afterBlank = $d.Paragraphs.Item(6)
$insPos = $afterBlank.Range.Start
$insRange2 = $d.Range($insPos, $insPos)
$insRange2.InsertParagraphBefore()
$insRange2.InsertParagraphBefore()
$insRange2.InsertParagraphBefore()

$d.Paragraphs.Item(6).Range.Text = "Biodiversity"
$d.Paragraphs.Item(8).Range.Text = "PI: Carla Hurt "

# ---------------------------------------------------------------------
# 3) Remove the scratch / brainstorming paragraphs near the end of the
#    document (from "Start writing down ..." through "Talk about
#    tangible applications achieved with sequence ").
# ---------------------------------------------------------------------
$scratchStartPara = $null
$scratchEndPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $txt = $d.Paragraphs.Item($i).Range.Text
    if ($txt -like "Start writing down*") {
        $scratchStartPara = $i
    }
    if ($txt -like "Talk about tangible applications achieved with sequence*") {
        $scratchEndPara = $i
    }
}

if ($scratchStartPara -ne $null -and $scratchEndPara -ne $null) {
    $scratchStart = $d.Paragraphs.Item($scratchStartPara).Range.Start
    $scratchEnd = $d.Paragraphs.Item($scratchEndPara).Range.End
    $d.Range($scratchStart, $scratchEnd).Delete()
}

Write-Output "Done. Paragraph count now: $($d.Paragraphs.Count)"
